$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.972.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.526.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.73"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.68%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.138.29"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +12.07%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.932.50"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.528.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.78"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "399.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.88"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.547"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.15"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.16"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.15"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.28"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.891.42"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0745"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.55"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "351.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.40%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.08"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.86"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.68%  "
